# Selecionar Receita — diagram rework
# - merges the "PesquisarReceita include" + "Selecionar Receita" + "Buscar Receita
#   Selecionada" normal-flow into two renumbered steps, drops the "Estar Logado no
#   Sistema" pre-condition in favour of "Pesquisar Receita", and renumbers the
#   cancellation exception flow from 2.x to 1.x.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 9 (empty C9 / "3. Buscar Receita Selecionada" in D9) is being removed;
# row 8 becomes the new last row of that B6:B8 block, so first clone row 9's
# (bottom-bordered) formatting onto row 8 before deleting row 9, keeping the
# borders/merges consistent with the rest of the table.
$ws.Range("C9:D9").Copy() | Out-Null
$ws.Range("C8:D8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Remove the now-redundant row; rows 10-12 shift up to become rows 9-11.
$ws.Rows.Item(9).Delete()

# Pré condição: "Estar Logado no Sistema" -> "Pesquisar Receita"
$ws.Range("C4").Value = "Pesquisar Receita"

# Cenário Normal, passo 1: "1. «include» PesquisarReceita" -> "1. Selecionar Receita"
$ws.Range("C7").Value = "1. Selecionar Receita"

# Cenário Normal, passo 2 now lives in column D of row 8 (renumbered 3 -> 2);
# clear the old column C text that used to hold "2. Selecionar Receita".
$ws.Range("C8").ClearContents()
$ws.Range("D8").Value = "2. Buscar Receita Selecionada"

# Cenário Excessão steps renumbered from 2.x to 1.x
$ws.Range("C9").Value = "1.1. Informa que quer cancelar seleção"
$ws.Range("D10").Value = "1.2. Cancela Selação de Receita"

# Restore the recorded selection/active cell for the sheet view
$ws.Range("C4:D4").Select()
